$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: split the run " 现在，风国大将军知道了敌军共有" into two runs:
#         " " and "现在，风国大将军知道了敌军共有"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(" 现在，风国大将军知道了敌军共有")
if ($found) {
    $rng.Text = ""
    $rng.InsertBefore("现在，风国大将军知道了敌军共有")
    $splitPoint = $d.Range($rng.Start, $rng.Start)
    $splitPoint.InsertBefore(" ")
}

# ---------------------------------------------------------------------
# Hunk 2: add the "不能构成环" clause about the road network
# ---------------------------------------------------------------------
$rng = $d.Content
$old2 = "号城市是他们的首都。在这些城市之间，有一些单向道路，并且保证从首都可以到达其他所有城市。大将军获得很多情报，每条情报表示敌军会从首都向一些城市增兵，大将军希望知道，有多少个城市是所有增兵的必经之地。（敌军如果会派遣"
$new2 = "号城市是他们的首都。在这些城市之间，有一些单向道路，并且保证从首都可以到达其他所有城市，且这些单向道路不能构成环（即不会出现从一个地方出发，走着走着又回到那个地方）。大将军获得很多情报，每条情报表示敌军会从首都向一些城市增兵，大将军希望知道，有多少个城市是所有增兵的必经之地。（敌军如果会派遣"
$found = $rng.Find.Execute($old2)
if ($found) {
    $rng.Text = $new2
}

# ---------------------------------------------------------------------
# Hunk 3: clarify that multiple operations may be chosen each time
# ---------------------------------------------------------------------
$rng = $d.Content
$old3 = "（也可以不操作），则称"
$new3 = "（可以进行多次操作，每次从两种操作中选择一个进行，也可以不操作），则称"
$found = $rng.Find.Execute($old3)
if ($found) {
    $rng.Text = $new3
}

# ---------------------------------------------------------------------
# Hunk 4: "字符串" -> "组成的字符串"
# ---------------------------------------------------------------------
$rng = $d.Content
$old4 = "夏荷给了冬雪两个由'A'和'B'字符串："
$new4 = "夏荷给了冬雪两个由'A'和'B'组成的字符串："
$found = $rng.Find.Execute($old4)
if ($found) {
    $rng.Text = $new4
}

# ---------------------------------------------------------------------
# Hunk 5: merge the " " run and "昆阳在给夏荷出题，题目是这样的：" run
#         into a single run
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("昆阳在给夏荷出题，题目是这样的：")
if ($found) {
    $mergeRng = $d.Range($rng.Start - 1, $rng.End)
    # force a real content change first so the two runs actually merge
    # into one (setting identical text is treated as a no-op)
    $mergeRng.Text = " 昆阳在给夏荷出题，题目是这样的：`u{0001}"
    $rng2 = $d.Content
    $found2 = $rng2.Find.Execute(" 昆阳在给夏荷出题，题目是这样的：`u{0001}")
    if ($found2) {
        $rng2.Text = " 昆阳在给夏荷出题，题目是这样的："
    }
}

# ---------------------------------------------------------------------
# Hunk 6: drop the leading " " run before "给你一棵包含"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("给你一棵包含")
if ($found) {
    $mergeRng = $d.Range($rng.Start - 1, $rng.End)
    $mergeRng.Text = "给你一棵包含"
}

# ---------------------------------------------------------------------
# Hunk 10: "1号店" -> "1号点" in the explanation sentence
# ---------------------------------------------------------------------
$rng = $d.Content
$old10 = "，表示先到1号店，再到2号点，然后重复任意次1、 2（可以是0次），最后到达3.对于"
$new10 = "，表示先到1号点，再到2号点，然后重复任意次1、 2（可以是0次），最后到达3.对于"
$found = $rng.Find.Execute($old10)
if ($found) {
    $rng.Text = $new10
}

Write-Output "done"
